$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 686, shifting existing rows 686:713 down to 687:714
$ws.Rows.Item(686).Insert()

# Populate the newly inserted row 686 with the new record
$ws.Range("A686").Value = 10
$ws.Range("B686").Value = "Vega Modelo de Temuco"
$ws.Range("C686").Value = "La Araucanía"
$ws.Range("D686").Value = 44747
$ws.Range("E686").Value = 9
$ws.Range("F686").Value = 100112006
$ws.Range("G686").Value = "Repollo"
$ws.Range("H686").Value = "Crespo record"
$ws.Range("I686").Value = "Primera"
$ws.Range("J686").Value = 750
$ws.Range("K686").Value = 1300
$ws.Range("L686").Value = 1300
$ws.Range("M686").Value = 1300
$ws.Range("N686").Value = "$/unidad"
$ws.Range("O686").Value = "Región del Maule"
$ws.Range("P686").Value = 1300
$ws.Range("Q686").Value = 1
$ws.Range("R686").Value = "Hortaliza"
